$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44306
$ws.Range("M4").Value = 80

# Row 5
$ws.Range("D5").Value = 44323
$ws.Range("M5").Value = 80

# Row 6
$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 120

# Row 7
$ws.Range("D7").Value = 44322
$ws.Range("M7").Value = 60

# Row 8
$ws.Range("D8").Value = 44313
$ws.Range("M8").Value = 120

# Row 9
$ws.Range("D9").Value = 44302

# Row 10
$ws.Range("D10").Value = 44330
$ws.Range("M10").Value = 60
